$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newest Argent (silver) price datapoint. The sheet's date/value
# columns store plain text (matching the rows above), so write the new row
# as formulas that yield the literal text, then paste-special as values so
# the cells end up holding static text (not numbers/dates and not live
# formulas) with the sheet's default (unformatted) cell style - exactly
# like rows 63-80 above it.
$ws.Range("A81").Formula = "=""2025-01-18"""
$ws.Range("B81").Formula = "=""42.6"""
$ws.Range("A81:B81").Copy()
$ws.Range("A81:B81").PasteSpecial(-4163)
$excel.CutCopyMode = 0

